# Executive presentation trim-down:
# Keep only slide 1 (title), slide 7 (Why This Solution?), slide 8
# (Business Value - Financial Impact) and slide 11 (Risk Mitigation),
# which become the new slides 1-4. All other slides are removed.
# Additionally, remove stray bold formatting from a few table cells.

$p = $ppt.ActivePresentation

# Delete slides from highest index to lowest so earlier indices are stable.
$slidesToDelete = @(17,16,15,14,13,12,10,9,6,5,4,3,2)
foreach ($idx in $slidesToDelete) {
    $p.Slides.Item($idx).Delete()
}

Write-Output "Slides remaining: $($p.Slides.Count)"

# --- New slide 2 (formerly slide 7): "Why This Solution?" table ---
# Row 4 ("[Current limitation 3]" / "[Our advantage 3]") loses its bold.
$s2 = $p.Slides.Item(2)
$tbl2 = $s2.Shapes.Item(3).Table
$tbl2.Cell(4,1).Shape.TextFrame.TextRange.Font.Bold = 0
$tbl2.Cell(4,2).Shape.TextFrame.TextRange.Font.Bold = 0

# --- New slide 3 (formerly slide 8): "Business Value - Financial Impact" ---
# Header row ("Metric" / "Value") and the ROI row lose their bold.
$s3 = $p.Slides.Item(3)
$tbl3 = $s3.Shapes.Item(3).Table
$tbl3.Cell(1,1).Shape.TextFrame.TextRange.Font.Bold = 0
$tbl3.Cell(1,2).Shape.TextFrame.TextRange.Font.Bold = 0
$tbl3.Cell(6,1).Shape.TextFrame.TextRange.Font.Bold = 0
$tbl3.Cell(6,2).Shape.TextFrame.TextRange.Font.Bold = 0

# --- New slide 4 (formerly slide 11): "Risk Mitigation" ---
# Header row ("Risk" / "Mitigation Strategy" / "Success Probability")
# and the "[Risk 3]" row lose their bold.
$s4 = $p.Slides.Item(4)
$tbl4 = $s4.Shapes.Item(3).Table
$tbl4.Cell(1,1).Shape.TextFrame.TextRange.Font.Bold = 0
$tbl4.Cell(1,2).Shape.TextFrame.TextRange.Font.Bold = 0
$tbl4.Cell(1,3).Shape.TextFrame.TextRange.Font.Bold = 0
$tbl4.Cell(4,1).Shape.TextFrame.TextRange.Font.Bold = 0
$tbl4.Cell(4,2).Shape.TextFrame.TextRange.Font.Bold = 0
$tbl4.Cell(4,3).Shape.TextFrame.TextRange.Font.Bold = 0

Write-Output "Done"
